$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 157; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 163; Resultado = "Acierto"; Profit = 0.67 },
    @{ Row = 165; Resultado = "Acierto"; Profit = 1.1 },
    @{ Row = 166; Resultado = "Acierto"; Profit = 0.33 },
    @{ Row = 167; Resultado = "Fallo";   Profit = -1 },
    @{ Row = 168; Resultado = "Acierto"; Profit = 0.67 },
    @{ Row = 170; Resultado = "Fallo";   Profit = -1 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 7).Value = $u.Resultado
    $ws.Cells.Item($u.Row, 8).Value = $u.Profit
}
